$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) holds text-formatted numbers (e.g. "69.503.62", "1.00").
# Force text format first so Excel does not coerce these into floating-point
# numbers (which would lose trailing zeros / thousands-dot formatting).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '69.503.62'
$ws.Range('E2').Value = '  -0.31%  '
$ws.Range('D3').Value = '3.743.63'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').Value = '612.41'
$ws.Range('E5').Value = '  -0.13%  '
$ws.Range('D6').Value = '177.42'
$ws.Range('E6').Value = '  -0.30%  '
$ws.Range('D7').Value = '3.743.94'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('D9').Value = '0.526'
$ws.Range('E9').Value = '  -2.66%  '
$ws.Range('D10').Value = '0.165'
$ws.Range('E10').Value = '  -0.91%  '
$ws.Range('D11').Value = '6.56'
$ws.Range('E11').Value = '  +3.15%  '
$ws.Range('D12').Value = '0.479'
$ws.Range('E12').Value = '  -4.09%  '
$ws.Range('D13').Value = '39.75'
$ws.Range('E13').Value = '  -2.87%  '
$ws.Range('D14').Value = '0.0000252'
$ws.Range('E14').Value = '  -0.92%  '
$ws.Range('D15').Value = '4.373.96'
$ws.Range('E15').Value = '  +0.19%  '
$ws.Range('D16').Value = '3.741.37'
$ws.Range('E16').Value = '  +0.04%  '
$ws.Range('D17').Value = '69.604.95'
$ws.Range('E17').Value = '  -0.27%  '
$ws.Range('E18').Value = '  -2.58%  '
$ws.Range('D19').Value = '7.41'
$ws.Range('E19').Value = '  -2.64%  '
$ws.Range('D20').Value = '16.32'
$ws.Range('E20').Value = '  -2.36%  '
$ws.Range('D21').Value = '499.26'
$ws.Range('E21').Value = '  -3.10%  '
$ws.Range('D22').Value = '9.13'
$ws.Range('E22').Value = '  -4.72%  '
$ws.Range('D23').Value = '0.717'
$ws.Range('E23').Value = '  -1.52%  '
$ws.Range('D24').Value = '2.56'
$ws.Range('E24').Value = '  +1.95%  '
$ws.Range('D25').Value = '85.75'
$ws.Range('E25').Value = '  -2.71%  '
$ws.Range('D26').Value = '12.80'
$ws.Range('E26').Value = '  -4.35%  '
$ws.Range('D27').Value = '10.96'
$ws.Range('E27').Value = '  -1.11%  '
$ws.Range('E28').Value = '  +4.81%  '
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('E30').Value = '  -2.51%  '
$ws.Range('E31').Value = '  +1.74%  '
$ws.Range('D32').Value = '7.98'
$ws.Range('E32').Value = '  +1.57%  '
$ws.Range('D33').Value = '30.29'
$ws.Range('E33').Value = '  -3.41%  '
$ws.Range('D34').Value = '0.112'
$ws.Range('E34').Value = '  -2.65%  '
$ws.Range('E35').Value = '  +0.15%  '
$ws.Range('E36').Value = '  +1.02%  '
$ws.Range('D37').Value = '6.07'
$ws.Range('E37').Value = '  -2.47%  '
$ws.Range('D38').Value = '0.345'
$ws.Range('E38').Value = '  +1.58%  '
$ws.Range('E39').Value = '  +2.76%  '
$ws.Range('D40').Value = '3.04'
$ws.Range('E40').Value = '  +11.22%  '
$ws.Range('D41').Value = '441.01'
$ws.Range('E41').Value = '  +3.79%  '
$ws.Range('E42').Value = '  -5.90%  '
$ws.Range('D43').Value = '49.69'
$ws.Range('E43').Value = '  -3.01%  '
$ws.Range('D44').Value = '44.23'
$ws.Range('E44').Value = '  -0.59%  '
$ws.Range('D45').Value = '8.50'
$ws.Range('E45').Value = '  -3.78%  '
$ws.Range('D46').Value = '2.943.08'
$ws.Range('E46').Value = '  -4.87%  '
$ws.Range('D47').Value = '0.0357'
$ws.Range('E47').Value = '  -2.26%  '
$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D48').Value = '1.00'
$ws.Range('E48').Value = '  -0.03%  '
$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').Value = '138.42'
$ws.Range('E49').Value = '  +1.90%  '
$ws.Range('D50').Value = '26.87'
$ws.Range('E50').Value = '  -3.50%  '
$ws.Range('D51').Value = '2.46'
$ws.Range('E51').Value = '  -2.50%  '
